$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("main")
$ws.Activate()

# Row 57: reword the "Expected Output" text (shared by rows 57-59) and fill
# in the previously-blank "Actual output if unexpected" / "Success?" /
# "Comments" columns.
$ws.Range("D57").Value = 'After this string is entered, it is copied to the "New destination string."'
$ws.Range("E57").Value = "SAME"
$ws.Range("F57").Value = "PASS"
$ws.Range("G57").Value = "Destination string is reset to empty as entered string is copied to a new destination string"

# Row 58: same "Actual output if unexpected" / "Success?" / "Comments" fill-in
# (Expected Output already held the same text, which was just reworded above).
$ws.Range("D58").Value = 'After this string is entered, it is copied to the "New destination string."'
$ws.Range("E58").Value = "SAME"
$ws.Range("F58").Value = "PASS"
$ws.Range("G58").Value = "Destination string is reset to empty as entered string is copied to a new destination string"

# Row 59: gains an "Expected Output" entry (same text) plus the same
# "Actual output" / "Success?" / "Comments" values, and its row grows to a
# two-line height to match rows 57-58.
$ws.Range("D59").Value = 'After this string is entered, it is copied to the "New destination string."'
$ws.Range("E59").Value = "SAME"
$ws.Range("F59").Value = "PASS"
$ws.Range("G59").Value = "Destination string is reset to empty as entered string is copied to a new destination string"
$ws.Rows.Item(59).RowHeight = 26.4

# Update the saved view state: scrolled so row 55 / column B is the
# top-left visible cell, with H59 as the active selection.
$win = $excel.ActiveWindow
$win.ScrollRow = 55
$win.ScrollColumn = 2
$ws.Range("H59").Select()
